# Apply the cyclic row-data shift for rows 2-6 of sheet "Artfynd":
# each row's full record moves down one row, and the last
# row (6) wraps its data around to the top (row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 64987400
$ws.Range("B2").Value = 89338
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 'Stjärntagging'
$ws.Range("G2").Value = 'Asterodon ferruginosus'
$ws.Range("H2").Value = 'Pat.'
$ws.Range("Q2").Value = 622582.8630347433
$ws.Range("R2").Value = 7259002.777051079
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '2014-08-24'
$ws.Range("Y2").ClearFormats()
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = '2014-08-24'
$ws.Range("AA2").ClearFormats()
$ws.Range("AX2").Value = 'Patrik Nygren'
$ws.Range("A3").Value = 67978160
$ws.Range("B3").Value = 89317
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 3242
$ws.Range("F3").Value = 'Vitplätt'
$ws.Range("G3").Value = 'Chaetodermella luna'
$ws.Range("H3").Value = '(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert'
$ws.Range("Q3").Value = 622342.7052900216
$ws.Range("R3").Value = 7259347.332542222
$ws.Range("A4").Value = 67978157
$ws.Range("B4").Value = 89633
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 'Fläckporing'
$ws.Range("G4").Value = 'Anthoporia albobrunnea'
$ws.Range("H4").Value = '(Romell) Karasiński & Niemelä'
$ws.Range("J4").Value = 'fruktkroppar'
$ws.Range("P4").Value = 'Abmoberget, Ly lm'
$ws.Range("Q4").Value = 622438.9085855351
$ws.Range("R4").Value = 7259490.537392656
$ws.Range("S4").Value = 10
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = '2017-10-05'
$ws.Range("Y4").ClearFormats()
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = '2017-10-05'
$ws.Range("AA4").ClearFormats()
$ws.Range("AC4").Value = ""
$ws.Range("AW4").Value = 'Patrik Nygren'
$ws.Range("AX4").Value = 'Patrik Nygren, per-erik mukka'
$ws.Range("A5").Value = 69884788
$ws.Range("B5").Value = 89317
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 3242
$ws.Range("F5").Value = 'Vitplätt'
$ws.Range("G5").Value = 'Chaetodermella luna'
$ws.Range("H5").Value = '(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert'
$ws.Range("Q5").Value = 622346.1763020725
$ws.Range("R5").Value = 7259344.144440647
$ws.Range("A6").Value = 69886350
$ws.Range("B6").Value = 89633
$ws.Range("D6").Value = 'VU'
$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 'Fläckporing'
$ws.Range("G6").Value = 'Anthoporia albobrunnea'
$ws.Range("H6").Value = '(Romell) Karasiński & Niemelä'
$ws.Range("J6").Value = ""
$ws.Range("P6").Value = 'Ned Saxnäs, Ly lm'
$ws.Range("Q6").Value = 622435.978354881
$ws.Range("R6").Value = 7259480.831004807
$ws.Range("S6").Value = 5
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = '2017-10-04'
$ws.Range("Y6").ClearFormats()
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = '2017-10-04'
$ws.Range("AA6").ClearFormats()
$ws.Range("AC6").Value = 'Påträffad under Sveaskogs naturvärdesinventering'
$ws.Range("AW6").Value = 'Mimmi Persson'
$ws.Range("AX6").Value = 'Mimmi Persson'
